# Add a new "Environment" asset (plus a couple of related queue/message
# assets) to the Assets sheet, to facilitate switching environments,
# assets and queues.

$wb = $excel.ActiveWorkbook

# Before switching sheets, leave the Settings sheet's selection at A18
# (matches the cell that was last selected there).
$settings = $wb.Worksheets.Item("Settings")
$settings.Range("A18").Select()

# Switch to / activate the Assets sheet and append the new asset rows.
$ws = $wb.Worksheets.Item("Assets")
$ws.Activate()

$ws.Range("A2").Value = "Environment"
$ws.Range("B2").Value = "RFW-Environment"

$ws.Range("A3").Value = "Message"
$ws.Range("B3").Value = "RFW-TestMessage"

$ws.Range("A4").Value = "TransactionQueue"
$ws.Range("B4").Value = "RFW-TransactionQueue"

# Leave selection on the Assets sheet at B5, right below the new data.
$ws.Range("B5").Select()
